# refactor and review smoke test
# Insert a new "LoginData" worksheet right before "accounts" (i.e. right
# after "consumerAccounts"), populate it with the Admin/Auditor/Client
# smoke-test login matrix, wire up mailto hyperlinks, and make it the
# active sheet (which also drops the old tabSelected flag from "Login").

$wb = $excel.ActiveWorkbook

# --- Create & place the new worksheet ------------------------------------
$beforeSheet = $wb.Worksheets.Item("accounts")
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "LoginData"

# --- Cell values (write in this exact order so the shared-string table
#     lines up the same way Excel would have built it) -------------------
$ws.Range("B1").Value = "Admin"
$ws.Range("C1").Value = "Auditor"
$ws.Range("D1").Value = "Client"
$ws.Range("B2").Value = "admin.user01@gmail.com"
$ws.Range("C2").Value = "auditor.user01@gmail.com"
$ws.Range("D2").Value = "client.user01@gmail.com"
$ws.Range("A2").Value = "Valid User"
$ws.Range("A3").Value = "Valid User2"
$ws.Range("E1").Value = "Auditor Lead"
$ws.Range("B3").Value = "admin.user02@gmail.com"
$ws.Range("C3").Value = "auditor.user02gmail.com"
$ws.Range("D3").Value = "client.user02@gmail.com"

# --- Header row + left-hand label column: blue fill + thin box border ---
$headerRange = $ws.Range("B1:E1")
$headerRange.Interior.ThemeColor = 5
$headerRange.Borders.LineStyle = 1

$labelRange = $ws.Range("A2:A3")
$labelRange.Interior.ThemeColor = 5
$labelRange.Borders.LineStyle = 1

# --- Plain bordered (no fill) cells --------------------------------------
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("E2").Borders.LineStyle = 1
$ws.Range("E3").Borders.LineStyle = 1

# --- Hyperlinked email cells: Hyperlink style + thin box border ---------
$linkCells = @("B2", "C2", "D2", "B3", "C3", "D3")
foreach ($addr in $linkCells) {
    $cell = $ws.Range($addr)
    $cell.Style = "Hyperlink"
    $cell.Borders.LineStyle = 1
}

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:admin.user01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:auditor.user01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:client.user01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:admin.user02@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:auditor.user01@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:client.user02@gmail.com")

# --- Leftover formatted-but-empty cell (B6), plain Hyperlink style ------
$ws.Range("B6").Style = "Hyperlink"

# --- Column widths (approximate; engine applies its own px rounding) ----
$ws.Columns.Item(1).ColumnWidth = 13.39
$ws.Columns.Item(2).ColumnWidth = 23.53
$ws.Columns.Item(3).ColumnWidth = 24.39
$ws.Columns.Item(4).ColumnWidth = 26.67
$ws.Columns.Item(5).ColumnWidth = 23.67

# --- Selection / view -----------------------------------------------------
$ws.Range("C12").Select() | Out-Null
